# Apply edits to the "Notes" sheet:
#  - Update the Description text
#  - Update the Source text and add a new "Source-link" row right after it
#  - Update the license note and add a new "More information on licensing"
#    row right after it

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Notes")

# Insert a new blank row right after the current "Source:" row (row 4),
# shifting everything below down by one. The old row 5 (blank) is now row 6,
# "Notes:" moves from row 6 to row 7, and the old "It is provided..." row
# (originally row 13) is now row 14.
$ws.Rows("5").Insert()

# Insert a second new blank row right after the (now shifted) license row
# (originally row 13, now row 14), shifting everything below it down by one
# more. The license text row itself stays at row 14.
$ws.Rows("15").Insert()

# Update existing text in place.
$ws.Range("A2").Value = "Description: Average Dependency Ratio"
$ws.Range("A4").Value = "Source: Profiles of higher local governments 2014 - Uganda Bureau of Statistics."
$ws.Range("A14").Value = "It is licensed under a Creative Commons Attribution 4.0 International license."

# Fill in the two newly inserted rows.
$ws.Range("A5").Value = "Source-link: http://www.ubos.org/onlinefiles/uploads/ubos/2009_HLG_%20Abstract_printed/CIS+UPLOADS/Profiles%20of%20Higher%20Local%20Governments_June_2014.pdf"
$ws.Range("A15").Value = "More information on licensing is available here: https://creativecommons.org/licenses/by/4.0/"
